# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-23) was previously listed
# with periods in descending order (1909 -> 1902). It is now resorted in
# ascending order (1902 -> 1909), keeping each period's corresponding
# "Valor Mora" amount attached to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for column E (rows 16-23)
$ws.Range("E16").Value = "1902"
$ws.Range("E17").Value = "1903"
$ws.Range("E18").Value = "1904"
$ws.Range("E19").Value = "1905"
$ws.Range("E20").Value = "1906"
$ws.Range("E21").Value = "1907"
$ws.Range("E22").Value = "1908"
$ws.Range("E23").Value = "1909"

# Matching "Valor Mora" amounts follow their period to the new row position
$ws.Range("F16").Value = 8833
$ws.Range("F17").Value = 33125
$ws.Range("F18").Value = 33125
$ws.Range("F19").Value = 33125
$ws.Range("F20").Value = 33125
$ws.Range("F21").Value = 33125
$ws.Range("F22").Value = 33125
$ws.Range("F23").Value = 33125
